$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.79115355014801
$ws.Range("B1").Value = 4.305760383605957
$ws.Range("C1").Value = 1.720344424247742
$ws.Range("D1").Value = 0.8724375367164612
$ws.Range("E1").Value = 0.4707854688167572
